$d = $word.ActiveDocument

# Phase 1: replace each old value with a unique placeholder token to avoid
# collisions where a new value contains another entry old value as a substring.
$d.Content.Find.Execute("1+36=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN0>>", 2) | Out-Null
$d.Content.Find.Execute("79-0=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN1>>", 2) | Out-Null
$d.Content.Find.Execute("29-21=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN2>>", 2) | Out-Null
$d.Content.Find.Execute("81-48=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN3>>", 2) | Out-Null
$d.Content.Find.Execute("34+26=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN4>>", 2) | Out-Null
$d.Content.Find.Execute("53-29=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN5>>", 2) | Out-Null
$d.Content.Find.Execute("37+61=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN6>>", 2) | Out-Null
$d.Content.Find.Execute("96-54=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN7>>", 2) | Out-Null
$d.Content.Find.Execute("7+46=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN8>>", 2) | Out-Null
$d.Content.Find.Execute("51-26=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN9>>", 2) | Out-Null
$d.Content.Find.Execute("85-57=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN10>>", 2) | Out-Null
$d.Content.Find.Execute("5+35=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN11>>", 2) | Out-Null
$d.Content.Find.Execute("24+29=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN12>>", 2) | Out-Null
$d.Content.Find.Execute("42-29=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN13>>", 2) | Out-Null
$d.Content.Find.Execute("23+40=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN14>>", 2) | Out-Null
$d.Content.Find.Execute("87-36=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN15>>", 2) | Out-Null
$d.Content.Find.Execute("23+22=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN16>>", 2) | Out-Null
$d.Content.Find.Execute("29+41=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN17>>", 2) | Out-Null
$d.Content.Find.Execute("82-64=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN18>>", 2) | Out-Null
$d.Content.Find.Execute("32-5=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN19>>", 2) | Out-Null
$d.Content.Find.Execute("69+28=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN20>>", 2) | Out-Null
$d.Content.Find.Execute("13-3=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN21>>", 2) | Out-Null
$d.Content.Find.Execute("59+39=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN22>>", 2) | Out-Null
$d.Content.Find.Execute("56+31=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN23>>", 2) | Out-Null
$d.Content.Find.Execute("76-33=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN24>>", 2) | Out-Null
$d.Content.Find.Execute("95-69=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN25>>", 2) | Out-Null
$d.Content.Find.Execute("6+66=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN26>>", 2) | Out-Null
$d.Content.Find.Execute("36-8=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN27>>", 2) | Out-Null
$d.Content.Find.Execute("59-18=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN28>>", 2) | Out-Null
$d.Content.Find.Execute("61-44=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN29>>", 2) | Out-Null
$d.Content.Find.Execute("81-23=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN30>>", 2) | Out-Null
$d.Content.Find.Execute("31+4=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN31>>", 2) | Out-Null
$d.Content.Find.Execute("56-16=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN32>>", 2) | Out-Null
$d.Content.Find.Execute("64+3=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN33>>", 2) | Out-Null
$d.Content.Find.Execute("83-21=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN34>>", 2) | Out-Null
$d.Content.Find.Execute("2+22=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN35>>", 2) | Out-Null
$d.Content.Find.Execute("35+22=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN36>>", 2) | Out-Null
$d.Content.Find.Execute("31-20=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN37>>", 2) | Out-Null
$d.Content.Find.Execute("78-76=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN38>>", 2) | Out-Null
$d.Content.Find.Execute("86+10=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN39>>", 2) | Out-Null
$d.Content.Find.Execute("96-91=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN40>>", 2) | Out-Null
$d.Content.Find.Execute("67-61=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN41>>", 2) | Out-Null
$d.Content.Find.Execute("43+22=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN42>>", 2) | Out-Null
$d.Content.Find.Execute("40-23=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN43>>", 2) | Out-Null
$d.Content.Find.Execute("21+45=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN44>>", 2) | Out-Null
$d.Content.Find.Execute("18+71=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN45>>", 2) | Out-Null
$d.Content.Find.Execute("58-0=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN46>>", 2) | Out-Null
$d.Content.Find.Execute("66+16=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN47>>", 2) | Out-Null
$d.Content.Find.Execute("10+3=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN48>>", 2) | Out-Null
$d.Content.Find.Execute("44+28=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN49>>", 2) | Out-Null
$d.Content.Find.Execute("42-35=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN50>>", 2) | Out-Null
$d.Content.Find.Execute("57-27=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN51>>", 2) | Out-Null
$d.Content.Find.Execute("68-0=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN52>>", 2) | Out-Null
$d.Content.Find.Execute("76-25=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN53>>", 2) | Out-Null
$d.Content.Find.Execute("68-66=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN54>>", 2) | Out-Null
$d.Content.Find.Execute("89-20=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN55>>", 2) | Out-Null
$d.Content.Find.Execute("29+6=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN56>>", 2) | Out-Null
$d.Content.Find.Execute("53-49=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN57>>", 2) | Out-Null
$d.Content.Find.Execute("17+76=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN58>>", 2) | Out-Null
$d.Content.Find.Execute("59+1=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN59>>", 2) | Out-Null
$d.Content.Find.Execute("28+30=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN60>>", 2) | Out-Null
$d.Content.Find.Execute("18-2=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN61>>", 2) | Out-Null
$d.Content.Find.Execute("0+3=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN62>>", 2) | Out-Null
$d.Content.Find.Execute("16+17=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN63>>", 2) | Out-Null
$d.Content.Find.Execute("28-12=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN64>>", 2) | Out-Null
$d.Content.Find.Execute("42+40=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN65>>", 2) | Out-Null
$d.Content.Find.Execute("5+14=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN66>>", 2) | Out-Null
$d.Content.Find.Execute("55-45=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN67>>", 2) | Out-Null
$d.Content.Find.Execute("11+40=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN68>>", 2) | Out-Null
$d.Content.Find.Execute("44+48=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN69>>", 2) | Out-Null
$d.Content.Find.Execute("69-29=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN70>>", 2) | Out-Null
$d.Content.Find.Execute("59-47=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN71>>", 2) | Out-Null
$d.Content.Find.Execute("71+7=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN72>>", 2) | Out-Null
$d.Content.Find.Execute("98-76=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN73>>", 2) | Out-Null
$d.Content.Find.Execute("78-47=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN74>>", 2) | Out-Null
$d.Content.Find.Execute("38+19=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN75>>", 2) | Out-Null
$d.Content.Find.Execute("17+63=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN76>>", 2) | Out-Null
$d.Content.Find.Execute("75+23=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN77>>", 2) | Out-Null
$d.Content.Find.Execute("43-9=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN78>>", 2) | Out-Null
$d.Content.Find.Execute("53-7=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN79>>", 2) | Out-Null
$d.Content.Find.Execute("50+28=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN80>>", 2) | Out-Null
$d.Content.Find.Execute("67-62=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN81>>", 2) | Out-Null
$d.Content.Find.Execute("82-76=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN82>>", 2) | Out-Null
$d.Content.Find.Execute("6+69=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN83>>", 2) | Out-Null
$d.Content.Find.Execute("45-17=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN84>>", 2) | Out-Null
$d.Content.Find.Execute("85-52=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN85>>", 2) | Out-Null
$d.Content.Find.Execute("74+4=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN86>>", 2) | Out-Null
$d.Content.Find.Execute("36-27=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN87>>", 2) | Out-Null
$d.Content.Find.Execute("35-20=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN88>>", 2) | Out-Null
$d.Content.Find.Execute("71+4=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN89>>", 2) | Out-Null
$d.Content.Find.Execute("7+75=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN90>>", 2) | Out-Null
$d.Content.Find.Execute("38+49=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN91>>", 2) | Out-Null
$d.Content.Find.Execute("89-63=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN92>>", 2) | Out-Null
$d.Content.Find.Execute("79-51=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN93>>", 2) | Out-Null
$d.Content.Find.Execute("14+10=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN94>>", 2) | Out-Null
$d.Content.Find.Execute("80-73=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN95>>", 2) | Out-Null
$d.Content.Find.Execute("83-49=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN96>>", 2) | Out-Null
$d.Content.Find.Execute("54+26=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN97>>", 2) | Out-Null
$d.Content.Find.Execute("46-11=", $true, $false, $false, $false, $false, $true, 1, $false, "<<TOKEN98>>", 2) | Out-Null

# Phase 2: replace placeholder tokens with the final values.
$d.Content.Find.Execute("<<TOKEN0>>", $true, $false, $false, $false, $false, $true, 1, $false, "44+24=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN1>>", $true, $false, $false, $false, $false, $true, 1, $false, "54-16=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN2>>", $true, $false, $false, $false, $false, $true, 1, $false, "13-7=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN3>>", $true, $false, $false, $false, $false, $true, 1, $false, "52+35=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN4>>", $true, $false, $false, $false, $false, $true, 1, $false, "99-88=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN5>>", $true, $false, $false, $false, $false, $true, 1, $false, "53-46=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN6>>", $true, $false, $false, $false, $false, $true, 1, $false, "81-2=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN7>>", $true, $false, $false, $false, $false, $true, 1, $false, "74+10=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN8>>", $true, $false, $false, $false, $false, $true, 1, $false, "53+14=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN9>>", $true, $false, $false, $false, $false, $true, 1, $false, "30+4=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN10>>", $true, $false, $false, $false, $false, $true, 1, $false, "30+14=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN11>>", $true, $false, $false, $false, $false, $true, 1, $false, "9+61=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN12>>", $true, $false, $false, $false, $false, $true, 1, $false, "60-35=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN13>>", $true, $false, $false, $false, $false, $true, 1, $false, "54+11=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN14>>", $true, $false, $false, $false, $false, $true, 1, $false, "69-49=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN15>>", $true, $false, $false, $false, $false, $true, 1, $false, "77+13=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN16>>", $true, $false, $false, $false, $false, $true, 1, $false, "18-0=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN17>>", $true, $false, $false, $false, $false, $true, 1, $false, "82-58=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN18>>", $true, $false, $false, $false, $false, $true, 1, $false, "35-22=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN19>>", $true, $false, $false, $false, $false, $true, 1, $false, "83+11=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN20>>", $true, $false, $false, $false, $false, $true, 1, $false, "39+9=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN21>>", $true, $false, $false, $false, $false, $true, 1, $false, "12+69=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN22>>", $true, $false, $false, $false, $false, $true, 1, $false, "56-27=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN23>>", $true, $false, $false, $false, $false, $true, 1, $false, "52+10=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN24>>", $true, $false, $false, $false, $false, $true, 1, $false, "0+32=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN25>>", $true, $false, $false, $false, $false, $true, 1, $false, "57-37=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN26>>", $true, $false, $false, $false, $false, $true, 1, $false, "86-10=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN27>>", $true, $false, $false, $false, $false, $true, 1, $false, "78-3=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN28>>", $true, $false, $false, $false, $false, $true, 1, $false, "13+4=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN29>>", $true, $false, $false, $false, $false, $true, 1, $false, "70-65=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN30>>", $true, $false, $false, $false, $false, $true, 1, $false, "97-93=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN31>>", $true, $false, $false, $false, $false, $true, 1, $false, "63-34=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN32>>", $true, $false, $false, $false, $false, $true, 1, $false, "7+73=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN33>>", $true, $false, $false, $false, $false, $true, 1, $false, "61-32=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN34>>", $true, $false, $false, $false, $false, $true, 1, $false, "62-48=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN35>>", $true, $false, $false, $false, $false, $true, 1, $false, "75-48=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN36>>", $true, $false, $false, $false, $false, $true, 1, $false, "70+3=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN37>>", $true, $false, $false, $false, $false, $true, 1, $false, "99-33=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN38>>", $true, $false, $false, $false, $false, $true, 1, $false, "44+22=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN39>>", $true, $false, $false, $false, $false, $true, 1, $false, "68-19=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN40>>", $true, $false, $false, $false, $false, $true, 1, $false, "51-38=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN41>>", $true, $false, $false, $false, $false, $true, 1, $false, "32-32=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN42>>", $true, $false, $false, $false, $false, $true, 1, $false, "21-16=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN43>>", $true, $false, $false, $false, $false, $true, 1, $false, "96-51=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN44>>", $true, $false, $false, $false, $false, $true, 1, $false, "43-23=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN45>>", $true, $false, $false, $false, $false, $true, 1, $false, "78-35=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN46>>", $true, $false, $false, $false, $false, $true, 1, $false, "75-61=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN47>>", $true, $false, $false, $false, $false, $true, 1, $false, "25+37=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN48>>", $true, $false, $false, $false, $false, $true, 1, $false, "79-63=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN49>>", $true, $false, $false, $false, $false, $true, 1, $false, "25+72=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN50>>", $true, $false, $false, $false, $false, $true, 1, $false, "4+74=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN51>>", $true, $false, $false, $false, $false, $true, 1, $false, "72+24=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN52>>", $true, $false, $false, $false, $false, $true, 1, $false, "0+20=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN53>>", $true, $false, $false, $false, $false, $true, 1, $false, "2+44=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN54>>", $true, $false, $false, $false, $false, $true, 1, $false, "98-46=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN55>>", $true, $false, $false, $false, $false, $true, 1, $false, "67-8=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN56>>", $true, $false, $false, $false, $false, $true, 1, $false, "26+31=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN57>>", $true, $false, $false, $false, $false, $true, 1, $false, "40+52=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN58>>", $true, $false, $false, $false, $false, $true, 1, $false, "84-0=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN59>>", $true, $false, $false, $false, $false, $true, 1, $false, "41-27=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN60>>", $true, $false, $false, $false, $false, $true, 1, $false, "32+66=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN61>>", $true, $false, $false, $false, $false, $true, 1, $false, "66+8=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN62>>", $true, $false, $false, $false, $false, $true, 1, $false, "62+19=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN63>>", $true, $false, $false, $false, $false, $true, 1, $false, "9-7=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN64>>", $true, $false, $false, $false, $false, $true, 1, $false, "99-67=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN65>>", $true, $false, $false, $false, $false, $true, 1, $false, "26+20=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN66>>", $true, $false, $false, $false, $false, $true, 1, $false, "61+21=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN67>>", $true, $false, $false, $false, $false, $true, 1, $false, "73+14=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN68>>", $true, $false, $false, $false, $false, $true, 1, $false, "64-41=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN69>>", $true, $false, $false, $false, $false, $true, 1, $false, "53-46=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN70>>", $true, $false, $false, $false, $false, $true, 1, $false, "7+67=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN71>>", $true, $false, $false, $false, $false, $true, 1, $false, "99-5=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN72>>", $true, $false, $false, $false, $false, $true, 1, $false, "91-21=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN73>>", $true, $false, $false, $false, $false, $true, 1, $false, "64-23=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN74>>", $true, $false, $false, $false, $false, $true, 1, $false, "77-62=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN75>>", $true, $false, $false, $false, $false, $true, 1, $false, "28+18=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN76>>", $true, $false, $false, $false, $false, $true, 1, $false, "3+6=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN77>>", $true, $false, $false, $false, $false, $true, 1, $false, "64-4=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN78>>", $true, $false, $false, $false, $false, $true, 1, $false, "88-16=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN79>>", $true, $false, $false, $false, $false, $true, 1, $false, "85-46=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN80>>", $true, $false, $false, $false, $false, $true, 1, $false, "17+13=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN81>>", $true, $false, $false, $false, $false, $true, 1, $false, "26+47=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN82>>", $true, $false, $false, $false, $false, $true, 1, $false, "6+47=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN83>>", $true, $false, $false, $false, $false, $true, 1, $false, "39-10=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN84>>", $true, $false, $false, $false, $false, $true, 1, $false, "98-30=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN85>>", $true, $false, $false, $false, $false, $true, 1, $false, "53+24=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN86>>", $true, $false, $false, $false, $false, $true, 1, $false, "91-54=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN87>>", $true, $false, $false, $false, $false, $true, 1, $false, "80-29=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN88>>", $true, $false, $false, $false, $false, $true, 1, $false, "25-13=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN89>>", $true, $false, $false, $false, $false, $true, 1, $false, "46+46=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN90>>", $true, $false, $false, $false, $false, $true, 1, $false, "46+44=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN91>>", $true, $false, $false, $false, $false, $true, 1, $false, "92-91=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN92>>", $true, $false, $false, $false, $false, $true, 1, $false, "80-59=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN93>>", $true, $false, $false, $false, $false, $true, 1, $false, "63-12=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN94>>", $true, $false, $false, $false, $false, $true, 1, $false, "89-87=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN95>>", $true, $false, $false, $false, $false, $true, 1, $false, "24+24=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN96>>", $true, $false, $false, $false, $false, $true, 1, $false, "98-71=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN97>>", $true, $false, $false, $false, $false, $true, 1, $false, "67+16=", 2) | Out-Null
$d.Content.Find.Execute("<<TOKEN98>>", $true, $false, $false, $false, $false, $true, 1, $false, "94-71=", 2) | Out-Null
